# Add thêm nhân sự Nguyễn Hữu Quang
# Updates the "Lương" sheet with the recalculated payroll figures
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

$ws.Range("B1").Value2 = 8
$ws.Range("B2").Value2 = 31
$ws.Range("B3").Value2 = 16607142.85714286
$ws.Range("B14").Value2 = 11071428.57142857
$ws.Range("B25").Value2 = 16607142.85714286
$ws.Range("B34").Value2 = 12195142.85714286
$ws.Range("B35").Value2 = 11071428.57142857
$ws.Range("B36").Value2 = 16607142.85714286
$ws.Range("B37").Value2 = 39873714.28571429
